$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 684.7
$ws.Range("I33").Value = 728
$ws.Range("K33").Value = 728
$ws.Range("M33").Value = -499
$ws.Range("H34").Value = 1966.9
$ws.Range("I34").Value = 1966.9
$ws.Range("K34").Value = 1966.9
$ws.Range("M34").Value = -1763.9
$ws.Range("H36").Value = 1966.9
$ws.Range("I36").Value = 1966.9
$ws.Range("K36").Value = 1966.9
$ws.Range("M36").Value = -1251.9
$ws.Range("H86").Value = 7708.3335
$ws.Range("I86").Value = 4942.4
$ws.Range("J86").Value = 11165.75
$ws.Range("K86").Value = 4942.4
$ws.Range("L86").Value = 11165.75
$ws.Range("M86").Value = -3819.4
$ws.Range("N86").Value = -13411.75
$ws.Range("H89").Value = 7708.3335
$ws.Range("I89").Value = 4942.4
$ws.Range("J89").Value = 11165.75
$ws.Range("K89").Value = 24712
$ws.Range("L89").Value = 55828.75
$ws.Range("M89").Value = -19096
$ws.Range("N89").Value = -67060.75
$ws.Range("H111").Value = 8549339
$ws.Range("I111").Value = 15876337
$ws.Range("J111").Value = 1173.6666
$ws.Range("K111").Value = 47629011
$ws.Range("L111").Value = 3520.9998
$ws.Range("M111").Value = -47625944
$ws.Range("N111").Value = -9654.9998
$ws.Range("H116").Value = 6175.4
$ws.Range("I116").Value = 4995
$ws.Range("K116").Value = 4995
$ws.Range("M116").Value = -1553
$ws.Range("H138").Value = 2674.725
$ws.Range("I138").Value = 1239.6207
$ws.Range("J138").Value = 3490.7646
$ws.Range("K138").Value = 3718.8621
$ws.Range("L138").Value = 10472.2938
$ws.Range("M138").Value = 1421.1379
$ws.Range("N138").Value = -20752.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3536877
$ws.Range("I2").Value = 4041834.5
$ws.Range("K2").Value = 4041834.5
$ws.Range("M2").Value = -4041721.5
$ws.Range("H32").Value = 9771.855
$ws.Range("I32").Value = 5897.6978
$ws.Range("J32").Value = 16179.115
$ws.Range("K32").Value = 5897.6978
$ws.Range("L32").Value = 16179.115
$ws.Range("M32").Value = -5610.6978
$ws.Range("N32").Value = -16753.115
$ws.Range("H45").Value = 7996877
$ws.Range("I45").Value = 15985925
$ws.Range("K45").Value = 15985925
$ws.Range("M45").Value = -15985548
$ws.Range("H96").Value = 66348.5
$ws.Range("J96").Value = 66348.5
$ws.Range("L96").Value = 66348.5
$ws.Range("N96").Value = -71840.5
$ws.Range("H97").Value = 662715.8
$ws.Range("I97").Value = 954555.9
$ws.Range("K97").Value = 954555.9
$ws.Range("M97").Value = -954059.9
$ws.Range("H116").Value = 3536877
$ws.Range("I116").Value = 4041834.5
$ws.Range("K116").Value = 4041834.5
$ws.Range("M116").Value = -4039540.5
$ws.Range("H132").Value = 2936.2307
$ws.Range("J132").Value = 5049.5
$ws.Range("L132").Value = 15148.5
$ws.Range("N132").Value = -20208.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3536877
$ws.Range("I3").Value = 4041834.5
$ws.Range("K3").Value = 4041834.5
$ws.Range("M3").Value = -4041720.5
$ws.Range("H8").Value = 442.8
$ws.Range("I8").Value = 366.5
$ws.Range("J8").Value = 493.66666
$ws.Range("K8").Value = 366.5
$ws.Range("L8").Value = 493.66666
$ws.Range("M8").Value = -226.5
$ws.Range("N8").Value = -773.66666
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 40
$ws.Range("N11").ClearContents()
$ws.Range("H86").Value = 3132984.8
$ws.Range("I86").Value = 4005838.5
$ws.Range("J86").Value = 15649.857
$ws.Range("K86").Value = 4005838.5
$ws.Range("L86").Value = 15649.857
$ws.Range("M86").Value = -4004715.5
$ws.Range("N86").Value = -17895.857
$ws.Range("H89").Value = 3132984.8
$ws.Range("I89").Value = 4005838.5
$ws.Range("J89").Value = 15649.857
$ws.Range("K89").Value = 20029192.5
$ws.Range("L89").Value = 78249.285
$ws.Range("M89").Value = -20023576.5
$ws.Range("N89").Value = -89481.285
$ws.Range("H94").Value = 2328981.2
$ws.Range("I94").Value = 3126322
$ws.Range("J94").Value = 9444.637000000001
$ws.Range("K94").Value = 3126322
$ws.Range("L94").Value = 9444.637000000001
$ws.Range("M94").Value = -3125871
$ws.Range("N94").Value = -10346.637
$ws.Range("H99").Value = 12989784
$ws.Range("I99").Value = 71429570
$ws.Range("J99").Value = 3165.5557
$ws.Range("K99").Value = 71429570
$ws.Range("L99").Value = 3165.5557
$ws.Range("M99").Value = -71428072
$ws.Range("N99").Value = -6161.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 128327
$ws.Range("J52").Value = 128327
$ws.Range("L52").Value = 128327
$ws.Range("N52").Value = -128915
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 4872.75
$ws.Range("I99").Value = 4779.6665
$ws.Range("J99").Value = 4965.8335
$ws.Range("K99").Value = 4779.6665
$ws.Range("L99").Value = 4965.8335
$ws.Range("M99").Value = -3281.6665
$ws.Range("N99").Value = -7961.8335
$ws.Range("H126").Value = 4872.75
$ws.Range("I126").Value = 4779.6665
$ws.Range("J126").Value = 4965.8335
$ws.Range("K126").Value = 14338.9995
$ws.Range("L126").Value = 14897.5005
$ws.Range("M126").Value = -11868.9995
$ws.Range("N126").Value = -19837.5005
$ws.Range("H132").Value = 102246.84
$ws.Range("I132").Value = 93482
$ws.Range("J132").Value = 114298.5
$ws.Range("K132").Value = 280446
$ws.Range("L132").Value = 342895.5
$ws.Range("M132").Value = -277916
$ws.Range("N132").Value = -347955.5
$ws.Range("H139").Value = 87490
$ws.Range("J139").Value = 87490
$ws.Range("L139").Value = 87490
$ws.Range("N139").Value = -97770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 2363
$ws.Range("I102").Value = 2363
$ws.Range("K102").Value = 7089
$ws.Range("M102").Value = -4655
$ws.Range("H131").Value = 12258033
$ws.Range("J131").Value = 17548060
$ws.Range("L131").Value = 52644180
$ws.Range("N131").Value = -52654260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H97").Value = 2165672
$ws.Range("I97").Value = 2646599
$ws.Range("J97").Value = 1499.5
$ws.Range("K97").Value = 2646599
$ws.Range("L97").Value = 1499.5
$ws.Range("M97").Value = -2646103
$ws.Range("N97").Value = -2491.5
$ws.Range("H109").Value = 48378.4
$ws.Range("J109").Value = 48378.4
$ws.Range("L109").Value = 48378.4
$ws.Range("N109").Value = -50458.4
$ws.Range("H137").Value = 106249
$ws.Range("J137").Value = 106249
$ws.Range("L137").Value = 106249
$ws.Range("N137").Value = -116449

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1797.2354
$ws.Range("I55").Value = 2708
$ws.Range("J55").Value = 1300.4546
$ws.Range("K55").Value = 2708
$ws.Range("L55").Value = 1300.4546
$ws.Range("M55").Value = -2535
$ws.Range("N55").Value = -1646.4546
$ws.Range("H61").Value = 3712110.5
$ws.Range("I61").Value = 4123326.2
$ws.Range("K61").Value = 4123326.2
$ws.Range("M61").Value = -4123124.2
$ws.Range("H68").Value = 1833.3334
$ws.Range("J68").Value = 2250
$ws.Range("L68").Value = 2250
$ws.Range("N68").Value = -3748
$ws.Range("H71").Value = 1833.3334
$ws.Range("J71").Value = 2250
$ws.Range("L71").Value = 11250
$ws.Range("N71").Value = -18738
$ws.Range("H93").Value = 10424613
$ws.Range("I93").Value = 15153303
$ws.Range("K93").Value = 15153303
$ws.Range("M93").Value = -15152055
$ws.Range("H113").Value = 3712110.5
$ws.Range("I113").Value = 4123326.2
$ws.Range("K113").Value = 4123326.2
$ws.Range("M113").Value = -4121156.2
$ws.Range("H132").Value = 4189.9707
$ws.Range("J132").Value = 6511.222
$ws.Range("L132").Value = 19533.666
$ws.Range("N132").Value = -24593.666
$ws.Range("H136").Value = 65406.727
$ws.Range("J136").Value = 6252.8237
$ws.Range("L136").Value = 18758.4711
$ws.Range("N136").Value = -23858.4711

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7782.6343
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 8418.594999999999
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 8418.594999999999
$ws.Range("M62").Value = -1276
$ws.Range("N62").Value = -9666.594999999999
$ws.Range("H65").Value = 7782.6343
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 8418.594999999999
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 42092.975
$ws.Range("M65").Value = -6380
$ws.Range("N65").Value = -48332.975
$ws.Range("H95").Value = 35174.25
$ws.Range("J95").Value = 35174.25
$ws.Range("L95").Value = 35174.25
$ws.Range("N95").Value = -40666.25
$ws.Range("H107").Value = 40001124
$ws.Range("I107").Value = 52631960
$ws.Range("K107").Value = 157895880
$ws.Range("M107").Value = -157893960
$ws.Range("H132").Value = 26604478
$ws.Range("I132").Value = 32260330
$ws.Range("J132").Value = 1557125.6
$ws.Range("K132").Value = 96780990
$ws.Range("L132").Value = 4671376.800000001
$ws.Range("M132").Value = -96778460
$ws.Range("N132").Value = -4676436.800000001
